# Updates cryptos list values (Price / Volume(1h)) per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value is a plain decimal number must be forced to
# Text format first, otherwise Excel auto-converts the assignment to a Number and
# the original text formatting (e.g. trailing zeros, thousand-dot grouping) is lost.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "37.354.44"
$ws.Range("E2").Value = "  +3.00%  "
$ws.Range("D3").Value = "2.095.34"
$ws.Range("E3").Value = "  +4.61%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "251.02"
$ws.Range("E5").Value = "  +3.04%  "
$ws.Range("D6").Value = "0.665"
$ws.Range("E6").Value = "  +1.19%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "54.37"
$ws.Range("E8").Value = "  +23.09%  "
$ws.Range("D9").Value = "61.81"
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("E10").Value = "  +3.83%  "
$ws.Range("D11").Value = "0.0744"
$ws.Range("E11").Value = "  +4.71%  "
$ws.Range("E12").Value = "  +8.39%  "
$ws.Range("D13").Value = "15.22"
$ws.Range("E13").Value = "  +6.21%  "
$ws.Range("D14").Value = "2.402.51"
$ws.Range("D15").Value = "0.838"
$ws.Range("E15").Value = "  +4.81%  "
$ws.Range("D16").Value = "2.103.82"
$ws.Range("E16").Value = "  +5.30%  "
$ws.Range("D17").Value = "5.20"
$ws.Range("E17").Value = "  +7.06%  "
$ws.Range("D18").Value = "37.286.78"
$ws.Range("E18").Value = "  +2.99%  "
$ws.Range("D19").Value = "72.86"
$ws.Range("D20").Value = "14.71"
$ws.Range("E20").Value = "  +15.56%  "
$ws.Range("D21").Value = "0.0₃0849"
$ws.Range("E21").Value = "  +4.99%  "
$ws.Range("D22").Value = "241.26"
$ws.Range("E22").Value = "  +2.15%  "
$ws.Range("E23").Value = "  +7.46%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("E25").Value = "  +2.79%  "
$ws.Range("D26").Value = "172.36"
$ws.Range("E26").Value = "  +4.63%  "
$ws.Range("E27").Value = "  +8.34%  "
$ws.Range("D28").Value = "20.69"
$ws.Range("E28").Value = "  +5.82%  "
$ws.Range("D29").Value = "2.01"
$ws.Range("E29").Value = "  +4.38%  "
$ws.Range("E30").Value = "  +2.58%  "
$ws.Range("D31").Value = "23.73"
$ws.Range("E31").Value = "  +8.06%  "
$ws.Range("D32").Value = "1.07"
$ws.Range("E32").Value = "  +27.14%  "
$ws.Range("E33").Value = "  +4.65%  "
$ws.Range("D34").Value = "0.0616"
$ws.Range("E34").Value = "  +6.42%  "
$ws.Range("E35").Value = "  +9.62%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("E37").Value = "  +5.16%  "
$ws.Range("D38").Value = "1.86"
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").Value = "2.26"
$ws.Range("E39").Value = "  +8.30%  "
$ws.Range("E40").Value = "  +1.93%  "
$ws.Range("D41").Value = "18.38"
$ws.Range("E41").Value = "  +15.88%  "
$ws.Range("E42").Value = "  +6.31%  "
$ws.Range("E43").Value = "  +5.95%  "
$ws.Range("D44").Value = "98.97"
$ws.Range("E44").Value = "  +4.21%  "
$ws.Range("D45").Value = "0.0931"
$ws.Range("E45").Value = "  +14.61%  "
$ws.Range("E46").Value = "  +1.67%  "
$ws.Range("D47").Value = "4.09"
$ws.Range("E47").Value = "  +103.42%  "
$ws.Range("D48").Value = "1.322.95"
$ws.Range("E48").Value = "  +1.40%  "
$ws.Range("E49").Value = "  +7.01%  "
$ws.Range("E50").Value = "  +15.28%  "
$ws.Range("D51").Value = "2.33"
$ws.Range("E51").Value = "  +7.94%  "

# Restore the default (unstyled) cell style now that the values are stored as text,
# so the only observable change is the cell content, matching the source data update.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D51").Style = "Normal"

Write-Output "Applied cryptos update"
